$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7 - this pushes old rows 7..33 down to 8..34,
# carrying their values/formatting with them (Excel's normal Insert behaviour).
$ws.Rows(7).Insert()

# Populate the newly inserted row 7 with the new weekly record.
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(7, 4).Value = [DateTime]"2023-04-19"
$ws.Cells.Item(7, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7, 5).Value = 15
$ws.Cells.Item(7, 6).Value = 100112013
$ws.Cells.Item(7, 7).Value = "Alcachofa"
$ws.Cells.Item(7, 8).Value = "Madrigal"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 160
$ws.Cells.Item(7, 11).Value = 23000
$ws.Cells.Item(7, 12).Value = 24000
$ws.Cells.Item(7, 13).Value = 23562
$ws.Cells.Item(7, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(7, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(7, 16).Value = 589
$ws.Cells.Item(7, 17).Value = 40
$ws.Cells.Item(7, 18).Value = "Hortaliza"
